$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 2
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 1

# Fill in row 7 values (A7:H7)
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0

# Fill in row 8 values (A8:H8)
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

# Fill in row 9 values (A9:H9) - all zero
$ws.Range("A9").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

# Update selection to I9 as in the diff
$ws.Range("I9").Select()

$wb.Save()
